$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells we touch stay text (avoid Excel auto-number coercion
# stripping meaningful trailing zeros, e.g. "43.30" -> 43.3).
$dCells = @("D2","D3","D5","D7","D8","D9","D12","D13","D14","D15","D16","D19","D21","D23","D24","D26","D27","D28","D29","D30","D31","D33","D34","D35","D38","D43","D44","D45","D47","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.487.39"
$ws.Range("E2").Value = "  +2.97%  "

$ws.Range("D3").Value = "1.604.54"
$ws.Range("E3").Value = "  +2.59%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "212.37"
$ws.Range("E5").Value = "  +0.96%  "

$ws.Range("E6").Value = "  +6.64%  "

$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").Value = "26.55"
$ws.Range("E8").Value = "  +5.42%  "

$ws.Range("D9").Value = "43.30"
$ws.Range("E9").Value = "  -1.40%  "

$ws.Range("E10").Value = "  +2.40%  "

$ws.Range("E11").Value = "  +2.37%  "

$ws.Range("D12").Value = "0.0910"
$ws.Range("E12").Value = "  +1.76%  "

$ws.Range("D13").Value = "1.836.00"
$ws.Range("E13").Value = "  +2.67%  "

$ws.Range("D14").Value = "1.614.76"
$ws.Range("E14").Value = "  +3.25%  "

$ws.Range("D15").Value = "29.522.40"
$ws.Range("E15").Value = "  +3.00%  "

$ws.Range("D16").Value = "0.533"
$ws.Range("E16").Value = "  +3.52%  "

$ws.Range("E17").Value = "  +1.26%  "

$ws.Range("E18").Value = "  +3.11%  "

$ws.Range("D19").Value = "242.17"
$ws.Range("E19").Value = "  +4.52%  "

$ws.Range("E20").Value = "  +3.63%  "

$ws.Range("D21").Value = "0.0₃0690"
$ws.Range("E21").Value = "  +2.18%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").Value = "3.98"
$ws.Range("E23").Value = "  +1.80%  "

$ws.Range("D24").Value = "9.14"
$ws.Range("E24").Value = "  +1.52%  "

$ws.Range("E25").Value = "  -1.02%  "

$ws.Range("D26").Value = "154.45"
$ws.Range("E26").Value = "  +2.47%  "

$ws.Range("D27").Value = "15.31"
$ws.Range("E27").Value = "  +3.53%  "

$ws.Range("D28").Value = "0.108"
$ws.Range("E28").Value = "  +5.07%  "

$ws.Range("D29").Value = "6.37"
$ws.Range("E29").Value = "  +2.39%  "

$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.16%  "

$ws.Range("D31").Value = "0.0472"
$ws.Range("E31").Value = "  +2.43%  "

$ws.Range("E32").Value = "  -0.23%  "

$ws.Range("D33").Value = "3.21"
$ws.Range("E33").Value = "  +1.63%  "

$ws.Range("D34").Value = "3.10"
$ws.Range("E34").Value = "  +3.86%  "

$ws.Range("D35").Value = "1.417.10"
$ws.Range("E35").Value = "  +1.90%  "

$ws.Range("E36").Value = "  -1.23%  "

$ws.Range("E37").Value = "  +2.96%  "

$ws.Range("D38").Value = "2.80"
$ws.Range("E38").Value = "  +5.56%  "

$ws.Range("E39").Value = "  +0.10%  "

$ws.Range("E40").Value = "  +1.90%  "

$ws.Range("E41").Value = "  +3.52%  "

$ws.Range("E42").Value = "  +0.31%  "

$ws.Range("D43").Value = "53.74"
$ws.Range("E43").Value = "  +24.16%  "

$ws.Range("D44").Value = "0.0485"
$ws.Range("E44").Value = "  +5.56%  "

$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  -0.18%  "

$ws.Range("E46").Value = "  +2.35%  "

$ws.Range("D47").Value = "65.57"
$ws.Range("E47").Value = "  +2.28%  "

$ws.Range("E48").Value = "  +0.12%  "

$ws.Range("D49").Value = "1.748.32"
$ws.Range("E49").Value = "  +2.47%  "

$ws.Range("D50").Value = "86.60"
$ws.Range("E50").Value = "  +1.42%  "

$ws.Range("D51").Value = "0.835"
$ws.Range("E51").Value = "  -3.83%  "
